$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 0.40867665887106
$ws.Range("B2").Value = 0.40867665887106
$ws.Range("C2").Value = 0.03279999135527759
$ws.Range("D2").Value = 0.3359704561717811
$ws.Range("E2").Value = 0.481382861570339
$ws.Range("F2").Value = 0.03433886264450848
$ws.Range("G2").Value = 0.03126112006604671
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 32
$ws.Range("L2").Value = 128
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 0.1
$ws.Range("O2").Value = 0.1
$ws.Range("Q2").Value = 0.0005
$ws.Range("R2").Value = 0.0001
$ws.Range("T2").Value = 3

# Row 3 updates
$ws.Range("A3").Value = 0.4030557613350829
$ws.Range("B3").Value = 0.40867665887106
$ws.Range("C3").Value = 0.03766601681321238
$ws.Range("D3").Value = 0.2984368542304953
$ws.Range("E3").Value = 0.5076746684396704
$ws.Range("F3").Value = 0.05460750591009855
$ws.Range("G3").Value = 0.02072452771632622
$ws.Range("H3").Value = 128
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 0.1
$ws.Range("O3").Value = 0.1
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0005
$ws.Range("S3").Value = 16
$ws.Range("T3").Value = 2
